# "support for Windows GUI autotest"
# Fill in the three new Windows-GUI-automation test steps (rows 6-8) on the
# "echo" sheet: step number, action (MENU/CLICK/INPUT), page, element and
# (for the last step) the expected output data - mirroring the existing
# GET/POST request rows above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: open NotePad, use the menu to reach "About Notepad" -----------
$ws.Range("B6").Value = "NotePad"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = "MENU"
$ws.Range("F6").Value = "无标题 - 记事本"
$ws.Range("G6").Value = "帮助(H)->关于记事本(A)"

# --- Row 7: click OK on the "About Notepad" dialog -------------------------
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = "CLICK"
$ws.Range("F7").Value = "关于`"记事本`""
$ws.Range("G7").Value = "确定"

# --- Row 8: type into the editor and verify the text -----------------------
$ws.Range("D8").Value = 3
$ws.Range("E8").Value = "INPUT"
$ws.Range("F8").Value = "无标题 - 记事本"
$ws.Range("G8").Value = "Edit"
$ws.Range("H8").Value = "hello world!"

# --- Column widths: widen the page/element/expected-result columns so the
#     longer Windows-GUI step descriptions are readable -------------------
$ws.Columns.Item(5).ColumnWidth = 10.57
$ws.Columns.Item(6).ColumnWidth = 21.71
$ws.Columns.Item(7).ColumnWidth = 24.71
$ws.Columns.Item(8).ColumnWidth = 65.71

# --- Move the active selection to the cell that was just filled in --------
$ws.Range("H8").Select()
